# Change the presentation's design/theme colour scheme from the
# "Integral" palette to the standard "Office Theme" palette (the
# built-in default PowerPoint theme). Equivalent to picking the
# "Office Theme" design from the Design tab in the UI.
#
# Per this host, themes are only editable through
# ThemeColorScheme.Colors(i).RGB / ThemeFontScheme.MajorFont/MinorFont
# (there is no supported "load a whole theme" operation), so each of
# the twelve theme colour slots is set individually to the target
# "Office" values. The font scheme and format scheme (fills / lines /
# effects) are already identical between the old and new theme, so
# only the colours need to change.

$p = $ppt.ActivePresentation
$tcs = $p.SlideMaster.Theme.ThemeColorScheme

function RGBVal([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

# Index order matches MsoThemeColorSchemeIndex:
# 1=Dark1 2=Light1 3=Dark2 4=Light2 5-10=Accent1-6 11=Hyperlink 12=FollowedHyperlink
$officeThemeColors = @(
    "000000",  # Dark 1
    "FFFFFF",  # Light 1
    "44546A",  # Dark 2
    "E7E6E6",  # Light 2
    "5B9BD5",  # Accent 1
    "ED7D31",  # Accent 2
    "A5A5A5",  # Accent 3
    "FFC000",  # Accent 4
    "4472C4",  # Accent 5
    "70AD47",  # Accent 6
    "0563C1",  # Hyperlink
    "954F72"   # Followed Hyperlink
)

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $tcs.Colors($i).RGB = RGBVal($officeThemeColors[$i - 1])
}
